$d = $word.ActiveDocument

# 1. Subtitle style should be based on Title (this also centers it, since
#    Title's paragraph formatting includes centered justification).
$subtitle = $d.Styles("Subtitle")
$subtitle.BaseStyle = $d.Styles("Title")

# 2. Remove the explicit (redundant) run color from the Subtitle style and
#    from the AbstractTitle style. The object model only lets us set colors
#    to a value (even "automatic" still serializes an explicit <w:color/>
#    element), so to truly drop the property we edit the canonical OOXML
#    package directly via WordOpenXML, which round-trips faithfully.
$xml = $d.WordOpenXML

$subtitleOld = '<w:rPr><w:rFonts w:eastAsiaTheme="majorEastAsia" w:cstheme="majorBidi"/><w:color w:val="595959" w:themeColor="text1" w:themeTint="A6"/><w:spacing w:val="15"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:style><w:style w:type="character" w:customStyle="1" w:styleId="SubtitleChar">'
$subtitleNew = '<w:rPr><w:rFonts w:eastAsiaTheme="majorEastAsia" w:cstheme="majorBidi"/><w:spacing w:val="15"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:style><w:style w:type="character" w:customStyle="1" w:styleId="SubtitleChar">'
$xml = $xml.Replace($subtitleOld, $subtitleNew)

$abstractTitleOld = '<w:rPr><w:b/><w:color w:val="345A8A"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr></w:style><w:style w:type="paragraph" w:customStyle="1" w:styleId="Abstract">'
$abstractTitleNew = '<w:rPr><w:b/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr></w:style><w:style w:type="paragraph" w:customStyle="1" w:styleId="Abstract">'
$xml = $xml.Replace($abstractTitleOld, $abstractTitleNew)

$d.WordOpenXML = $xml
